$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    8  = 4
    12 = -1
    14 = 0
    16 = -1
    17 = 0
    23 = -2
    24 = -1
    26 = 2
    29 = 0
    30 = -1
    36 = 0
    41 = -11
    43 = 4
    49 = -2
    51 = -12
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
